# Generate Report for Handoff
#
# The localization report is regenerated: the row that used to describe
# "2cb4c16e-...md" and the row that used to describe "94a090e5-...md" swap
# places (by content) on every sheet, the status moves from
# "Handed back: in sync with en-US" to "Ready for handoff", the handoff/
# handback timestamps advance, and the (now out of date) 2cb4c16e file
# picks up a new "stale handback" Error Detail message.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "94a090e5-d009-4903-bd37-f02204d0e7cd.md"
$ov.Range("B2").Value = "e2e\94a090e5-d009-4903-bd37-f02204d0e7cd.md"

$ov.Range("A3").Value = "2cb4c16e-5412-4c80-892d-9095861aed28.md"
$ov.Range("B3").Value = "e2e\2cb4c16e-5412-4c80-892d-9095861aed28.md"

$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-09-06 07:58:17"

# Hyperlinks on B2/B3 keep pointing at the same targets (rId2 -> 2cb4c16e,
# rId3 -> 94a090e5) but the row each target is displayed on swaps, so the
# display text needs to follow the new cell text.
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d2872f81f7108e0d0d04c7029c1e5b4076b38008/e2e/2cb4c16e-5412-4c80-892d-9095861aed28.md", "", "", "e2e\94a090e5-d009-4903-bd37-f02204d0e7cd.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d2872f81f7108e0d0d04c7029c1e5b4076b38008/e2e/94a090e5-d009-4903-bd37-f02204d0e7cd.md", "", "", "e2e\2cb4c16e-5412-4c80-892d-9095861aed28.md") | Out-Null
$ov.Range("B2:B3").Style = "Hyperlink"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "94a090e5-d009-4903-bd37-f02204d0e7cd.md"
$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("G2").Value = "94a090e5-d009-4903-bd37-f02204d0e7cd.758c38fd4887970a0b2ae490360ae7e80a095a14.zh-cn.xlf"
$zh.Range("I2").Value = "94a090e5-d009-4903-bd37-f02204d0e7cd.md"
$zh.Range("J2").Value = "94a090e5-d009-4903-bd37-f02204d0e7cd.758c38fd4887970a0b2ae490360ae7e80a095a14.zh-cn.xlf"

$zh.Range("A3").Value = "2cb4c16e-5412-4c80-892d-9095861aed28.md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("G3").Value = "2cb4c16e-5412-4c80-892d-9095861aed28.e59b79be0c545b7275be669867fd2e69b1d20aec.zh-cn.xlf"
$zh.Range("H3").Value = "2016-09-06 07:57:59"
$zh.Range("I3").Value = "2cb4c16e-5412-4c80-892d-9095861aed28.md"
$zh.Range("J3").Value = "2cb4c16e-5412-4c80-892d-9095861aed28.e59b79be0c545b7275be669867fd2e69b1d20aec.zh-cn.xlf"
$zh.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d2872f81f7108e0d0d04c7029c1e5b4076b38008/e2e/2cb4c16e-5412-4c80-892d-9095861aed28.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0f922970b0288dfc7a5d5a6ec209210c9adad90c/e2e/2cb4c16e-5412-4c80-892d-9095861aed28.md."

$zh.Range("A2:A2").ColumnWidth = 40
$zh.Columns.Item(16).ColumnWidth = 40

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d2872f81f7108e0d0d04c7029c1e5b4076b38008/e2e/94a090e5-d009-4903-bd37-f02204d0e7cd.md", "", "", "94a090e5-d009-4903-bd37-f02204d0e7cd.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/7a0b18b11b20601b3a84eba07aee1e754290d4ef/e2e/94a090e5-d009-4903-bd37-f02204d0e7cd.md", "", "", "94a090e5-d009-4903-bd37-f02204d0e7cd.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d2872f81f7108e0d0d04c7029c1e5b4076b38008/e2e/2cb4c16e-5412-4c80-892d-9095861aed28.md", "", "", "2cb4c16e-5412-4c80-892d-9095861aed28.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/7a0b18b11b20601b3a84eba07aee1e754290d4ef/e2e/2cb4c16e-5412-4c80-892d-9095861aed28.md", "", "", "2cb4c16e-5412-4c80-892d-9095861aed28.md") | Out-Null
$zh.Range("A2:A3").Style = "Hyperlink"
$zh.Range("I2:I3").Style = "Hyperlink"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "94a090e5-d009-4903-bd37-f02204d0e7cd.md"
$de.Range("C2").Value = "Ready for handoff"
$de.Range("G2").Value = "94a090e5-d009-4903-bd37-f02204d0e7cd.758c38fd4887970a0b2ae490360ae7e80a095a14.de-de.xlf"
$de.Range("H2").Value = "2016-09-06 07:58:17"
$de.Range("I2").Value = "94a090e5-d009-4903-bd37-f02204d0e7cd.md"
$de.Range("J2").Value = "94a090e5-d009-4903-bd37-f02204d0e7cd.758c38fd4887970a0b2ae490360ae7e80a095a14.de-de.xlf"

$de.Range("A3").Value = "2cb4c16e-5412-4c80-892d-9095861aed28.md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("G3").Value = "2cb4c16e-5412-4c80-892d-9095861aed28.e59b79be0c545b7275be669867fd2e69b1d20aec.de-de.xlf"
$de.Range("H3").Value = "2016-09-06 07:58:17"
$de.Range("I3").Value = "2cb4c16e-5412-4c80-892d-9095861aed28.md"
$de.Range("J3").Value = "2cb4c16e-5412-4c80-892d-9095861aed28.e59b79be0c545b7275be669867fd2e69b1d20aec.de-de.xlf"
$de.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d2872f81f7108e0d0d04c7029c1e5b4076b38008/e2e/2cb4c16e-5412-4c80-892d-9095861aed28.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0f922970b0288dfc7a5d5a6ec209210c9adad90c/e2e/2cb4c16e-5412-4c80-892d-9095861aed28.md."

$de.Columns.Item(16).ColumnWidth = 40

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d2872f81f7108e0d0d04c7029c1e5b4076b38008/e2e/94a090e5-d009-4903-bd37-f02204d0e7cd.md", "", "", "94a090e5-d009-4903-bd37-f02204d0e7cd.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/ac79d9466959a0d6e21cea6d373572c64d6fc582/e2e/94a090e5-d009-4903-bd37-f02204d0e7cd.md", "", "", "94a090e5-d009-4903-bd37-f02204d0e7cd.md") | Out-Null
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d2872f81f7108e0d0d04c7029c1e5b4076b38008/e2e/2cb4c16e-5412-4c80-892d-9095861aed28.md", "", "", "2cb4c16e-5412-4c80-892d-9095861aed28.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/ac79d9466959a0d6e21cea6d373572c64d6fc582/e2e/2cb4c16e-5412-4c80-892d-9095861aed28.md", "", "", "2cb4c16e-5412-4c80-892d-9095861aed28.md") | Out-Null
$de.Range("A2:A3").Style = "Hyperlink"
$de.Range("I2:I3").Style = "Hyperlink"
